# Update countries & provincias Spain
# Refresh of the COVID "Pais" dashboard data: new timestamp, updated
# daily counters for a number of countries, and four pairs of countries
# that swapped rank (so their names/labels trade places) as the
# underlying data source re-sorted between the 18:30 and 19:47 pulls.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = 'Datos actualizados a 23 de Julio de 2020 a las 19:47'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4129405
$ws.Range("C4").Value = 28530
$ws.Range("D4").Value = 1953958
$ws.Range("E4").Value = 2028782
$ws.Range("G4").Value = 482
$ws.Range("H4").Value = 146665

# India (row 6)
$ws.Range("B6").Value = 1286430
$ws.Range("C6").Value = 46746
$ws.Range("D6").Value = 815979
$ws.Range("E6").Value = 439808
$ws.Range("G6").Value = 753
$ws.Range("H6").Value = 30643

# Chile (row 11)
$ws.Range("B11").Value = 338759
$ws.Range("C11").Value = 2357
$ws.Range("D11").Value = 311431
$ws.Range("E11").Value = 18490
$ws.Range("G11").Value = 116
$ws.Range("H11").Value = 8838

# Turquia (row 18)
$ws.Range("B18").Value = 223315
$ws.Range("C18").Value = 913
$ws.Range("D18").Value = 206365
$ws.Range("E18").Value = 11387
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 5563

# Alemania (row 21)
$ws.Range("B21").Value = 204598
$ws.Range("C21").Value = 128
$ws.Range("E21").Value = 6415

# Francia (row 22)
$ws.Range("B22").Value = 179398
$ws.Range("C22").Value = 1062
$ws.Range("E22").Value = 69269

# Canada (row 24)
$ws.Range("B24").Value = 112485
$ws.Range("C24").Value = 245
$ws.Range("D24").Value = 98425
$ws.Range("E24").Value = 5190

# Israel (row 42)
$ws.Range("B42").Value = 57453
$ws.Range("C42").Value = 1368
$ws.Range("D42").Value = 23916
$ws.Range("E42").Value = 33097
$ws.Range("G42").Value = 10
$ws.Range("H42").Value = 440

# Irlanda (row 60)
$ws.Range("B60").Value = 25826
$ws.Range("C60").Value = 7
$ws.Range("E60").Value = 699
$ws.Range("G60").Value = 9
$ws.Range("H60").Value = 1763

# Argelia (row 61)
$ws.Range("B61").Value = 25484
$ws.Range("C61").Value = 612
$ws.Range("D61").Value = 17369
$ws.Range("E61").Value = 6991
$ws.Range("G61").Value = 13
$ws.Range("H61").Value = 1124

# Rows 66/67 swap rank: Nepal <-> Marruecos
$ws.Range("A66").Value = 'Marruecos'
$ws.Range("B66").Value = 18264
$ws.Range("C66").Value = 302
$ws.Range("D66").Value = 15872
$ws.Range("E66").Value = 2100
$ws.Range("G66").Value = 7
$ws.Range("H66").Value = 292

$ws.Range("A67").Value = 'Nepal'
$ws.Range("B67").Value = 18241
$ws.Range("C67").Value = 147
$ws.Range("D67").Value = 12840
$ws.Range("E67").Value = 5358
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 43

# Rows 108/109 swap rank: Somalia <-> Libano
$ws.Range("A108").Value = 'Libano'
$ws.Range("B108").Value = 3260
$ws.Range("C108").Value = 156
$ws.Range("D108").Value = 1619
$ws.Range("E108").Value = 1598
$ws.Range("H108").Value = 43

$ws.Range("A109").Value = 'Somalia'
$ws.Range("B109").Value = 3161
$ws.Range("D109").Value = 1495
$ws.Range("E109").Value = 1573
$ws.Range("H109").Value = 93

# Maldivas (row 110)
$ws.Range("B110").Value = 3120
$ws.Range("C110").Value = 17
$ws.Range("D110").Value = 2428
$ws.Range("E110").Value = 677

# Rows 112/113 swap rank: Congo <-> Mayotte
$ws.Range("A112").Value = 'Mayotte'
$ws.Range("B112").Value = 2862
$ws.Range("C112").Value = 23
$ws.Range("D112").Value = 2650
$ws.Range("E112").Value = 174
$ws.Range("H112").Value = 38

$ws.Range("A113").Value = 'Congo'
$ws.Range("B113").Value = 2851
$ws.Range("D113").Value = 666
$ws.Range("E113").Value = 2135
$ws.Range("H113").Value = 50

# Yemen (row 133)
$ws.Range("B133").Value = 1654
$ws.Range("C133").Value = 14
$ws.Range("D133").Value = 762
$ws.Range("G133").Value = 3
$ws.Range("H133").Value = 461

# Republica de Chipre (row 147)
$ws.Range("B147").Value = 1045
$ws.Range("C147").Value = 5
$ws.Range("E147").Value = 179

# Rows 210/211 swap rank: Groenlandia <-> Islas Malvinas (values identical)
$ws.Range("A210").Value = 'Islas Malvinas'
$ws.Range("A211").Value = 'Groenlandia'
